# Insert a new paragraph right after the "module-16: working with forms &
# user input" heading paragraph, containing the new subtitle line:
#   "module-16_1 starts here,, tutorial: 1-8"
# The new paragraph reuses the same paragraph/run formatting (sz/szCs 24,
# spacing after=0, line=276/auto) as its neighbours.

$d = $word.ActiveDocument

$headingText = "module-16: working with forms & user input"
$newText = "module-16_1 starts here,, tutorial: 1-8"

$heading = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $headingText) {
        $heading = $p
        break
    }
}

if ($heading -eq $null) {
    throw "Could not locate heading paragraph '$headingText'"
}

# InsertParagraphAfter() splits in a new empty paragraph right after the
# heading, copying the heading's paragraph formatting (pPr) and leaving an
# empty run whose rPr matches the paragraph mark's rPr.
$heading.Range.InsertParagraphAfter()

$newPara = $heading.Next()
$newPara.Range.Text = $newText

Write-Output "Inserted paragraph: $($newPara.Range.Text)"
